$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: update recipient / shipment details ---
$ws.Range("A2").Value = "SPOT4TONER Ι Κ Ε"
$ws.Range("B2").Value = "'800839540"
$ws.Range("D2").ClearContents()
$ws.Range("F2").Value = "'54621"
$ws.Range("G2").Value = "ΘΕΣΣΑΛΟΝΙΚΗ"
$ws.Range("H2").Value = "ΑΓΓΕΛΑΚΗ 3"
$ws.Range("N2").Value = 3
$ws.Range("O2").Value = "'45€"

# --- Row 3: new product line (blank shipment fields, new product) ---
$ws.Range("A3").Formula = '=""'
$ws.Range("B3").Formula = '=""'
$ws.Range("C3").Formula = '=""'
$ws.Range("D3").Formula = '=""'
$ws.Range("E3").Formula = '=""'
$ws.Range("F3").Formula = '=""'
$ws.Range("G3").Formula = '=""'
$ws.Range("H3").Formula = '=""'
$ws.Range("I3").Formula = '=""'
$ws.Range("J3").Formula = '=""'
$ws.Range("K3").Formula = '=""'
$ws.Range("L3").Value = "GPT-0535"
$ws.Range("M3").Value = "OKI TONER B432 BLACK ΣΥΜΒΑΤΟ 12000  ΣΕΛΙΔΕΣ"
$ws.Range("N3").Value = "'15"
$ws.Range("O3").Formula = '=""'
